$wb = $excel.ActiveWorkbook

$wsDish = $wb.Worksheets.Item("Dish")
$wsDish.Range("A3").Value = "dish0"
$wsDish.Range("A5").Select()

$wsOrder = $wb.Worksheets.Item("Order")
$wsOrder.Range("D1").Value = "ordertime"
$wsOrder.Range("A2").Value = "varchar2(20)"
$wsOrder.Range("B2").Value = "varchar2(20)"
$wsOrder.Range("C2").Value = "varchar2(10)"
$wsOrder.Range("D3").Select()

$wsCook = $wb.Worksheets.Item("Cook")
$wsCook.Range("C8").Select()

$wsCookFood = $wb.Worksheets.Item("CookFood")
$wsCookFood.Range("C1").Value = "cookfoodtime"
$wsCookFood.Range("D1").Value = "orderNo"
$wsCookFood.Range("E1").Value = "status"
$wsCookFood.Range("D2").Value = "varchar2(20)"
$wsCookFood.Range("E2").Value = "char(1)"
$wsCookFood.Range("D3").Value = "order0"
$wsCookFood.Range("E3").Value = "A"
$wsCookFood.Activate()
$wsCookFood.Range("D6").Select()
